$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.29347
$ws.Range("H2").Value = 0.8804099999999999
$ws.Range("I2").Value = 0.1501202107524681
$ws.Range("J2").Value = 0.1501202107524681
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.238740666666666
$ws.Range("N2").Value = 6.716221999999999
$ws.Range("O2").Value = 0.03262524687146927
$ws.Range("P2").Value = 0.03262524687146927
$ws.Range("Q2").Value = 0.6570032234466664
$ws.Range("R2").Value = 5.913029011019999
$ws.Range("S2").Value = 0.004897708936196269
$ws.Range("T2").Value = 0.004897708936196269
# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.29347
$ws.Range("H3").Value = 0.8804099999999999
$ws.Range("I3").Value = 0.1501202107524681
$ws.Range("J3").Value = 0.1501202107524681
$ws.Range("O3").Value = 0.8675125960695174
$ws.Range("P3").Value = 0.8675125960695174
$ws.Range("Q3").Value = 17.46986235058
$ws.Range("R3").Value = 157.22876115522
$ws.Range("S3").Value = 0.1302311737523767
$ws.Range("T3").Value = 0.1302311737523767
# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.29347
$ws.Range("H4").Value = 0.8804099999999999
$ws.Range("I4").Value = 0.1501202107524681
$ws.Range("J4").Value = 0.1501202107524681
$ws.Range("M4").Value = 6.852529666666666
$ws.Range("N4").Value = 20.557589
$ws.Range("O4").Value = 0.0998621570590134
$ws.Range("P4").Value = 0.0998621570590134
$ws.Range("Q4").Value = 2.011011881276666
$ws.Range("R4").Value = 18.09910693149
$ws.Range("S4").Value = 0.01499132806389517
$ws.Range("T4").Value = 0.01499132806389517
# Row 5
$ws.Range("I5").Value = 0.2180391153852712
$ws.Range("J5").Value = 0.2180391153852712
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 2.238740666666666
$ws.Range("N5").Value = 6.716221999999999
$ws.Range("O5").Value = 0.03262524687146927
$ws.Range("P5").Value = 0.03262524687146927
$ws.Range("Q5").Value = 0.9542512692164442
$ws.Range("R5").Value = 8.588261422947999
$ws.Range("S5").Value = 0.007113579967081247
$ws.Range("T5").Value = 0.007113579967081247
# Row 6
$ws.Range("I6").Value = 0.2180391153852712
$ws.Range("J6").Value = 0.2180391153852712
$ws.Range("O6").Value = 0.8675125960695174
$ws.Range("P6").Value = 0.8675125960695174
$ws.Range("S6").Value = 0.1891516790325777
$ws.Range("T6").Value = 0.1891516790325777
# Row 7
$ws.Range("I7").Value = 0.2180391153852712
$ws.Range("J7").Value = 0.2180391153852712
$ws.Range("M7").Value = 6.852529666666666
$ws.Range("N7").Value = 20.557589
$ws.Range("O7").Value = 0.0998621570590134
$ws.Range("P7").Value = 0.0998621570590134
$ws.Range("Q7").Value = 2.920854223591778
$ws.Range("R7").Value = 26.287688012326
$ws.Range("S7").Value = 0.0217738563856123
$ws.Range("T7").Value = 0.0217738563856123
# Row 8
$ws.Range("G8").Value = 1.235185333333333
$ws.Range("H8").Value = 3.705556
$ws.Range("I8").Value = 0.6318406738622607
$ws.Range("J8").Value = 0.6318406738622606
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 2.238740666666666
$ws.Range("N8").Value = 6.716221999999999
$ws.Range("O8").Value = 0.03262524687146927
$ws.Range("P8").Value = 0.03262524687146927
$ws.Range("Q8").Value = 2.765259636603555
$ws.Range("R8").Value = 24.887336729432
$ws.Range("S8").Value = 0.02061395796819175
$ws.Range("T8").Value = 0.02061395796819175
# Row 9
$ws.Range("G9").Value = 1.235185333333333
$ws.Range("H9").Value = 3.705556
$ws.Range("I9").Value = 0.6318406738622607
$ws.Range("J9").Value = 0.6318406738622606
$ws.Range("O9").Value = 0.8675125960695174
$ws.Range("P9").Value = 0.8675125960695174
$ws.Range("Q9").Value = 73.52887092646134
$ws.Range("R9").Value = 661.759838338152
$ws.Range("S9").Value = 0.5481297432845631
$ws.Range("T9").Value = 0.548129743284563
# Row 10
$ws.Range("G10").Value = 1.235185333333333
$ws.Range("H10").Value = 3.705556
$ws.Range("I10").Value = 0.6318406738622607
$ws.Range("J10").Value = 0.6318406738622606
$ws.Range("M10").Value = 6.852529666666666
$ws.Range("N10").Value = 20.557589
$ws.Range("O10").Value = 0.0998621570590134
$ws.Range("P10").Value = 0.0998621570590134
$ws.Range("Q10").Value = 8.464144140498222
$ws.Range("R10").Value = 76.177297264484
$ws.Range("S10").Value = 0.06309697260950595
$ws.Range("T10").Value = 0.06309697260950593
